# The underlying ranking data was regenerated with slightly different
# "matrices" scores, which changes the sort order of the table (rows are
# ranked by descending matrices score). Only the three rows around the
# middle of the table (rows 8-10, i.e. prior ranks 7-9) actually swap
# position; their per-person data (prolificid #, prolificid hash, name,
# race) moves together with them, while A (0-based index) and H (rank)
# stay sequential/fixed for the row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recalculated "matrices" scores for every row (F column).
$ws.Range("F2").Value  = 14.32124806351207
$ws.Range("F3").Value  = 13.25929927578149
$ws.Range("F4").Value  = 8.119568600285705
$ws.Range("F5").Value  = 7.411316600731239
$ws.Range("F6").Value  = 6.387209556654361
$ws.Range("F7").Value  = 6.159915716578424
$ws.Range("F8").Value  = 5.496872041548905
$ws.Range("F9").Value  = 5.269135301367183
$ws.Range("F10").Value = 5.262741384947466
$ws.Range("F11").Value = 3.214622224054206
$ws.Range("F12").Value = 1.200783564827204
$ws.Range("F13").Value = 0.09961791273931464

# Rows 8-10 (ranks 7-9) re-sort: Jamarii moves up to rank 7 (row 8),
# Matthew drops to rank 8 (row 9), Brennan drops to rank 9 (row 10).
# Each row keeps its own prolificid number/hash and race together.
$ws.Range("B8").Value = 32
$ws.Range("C8").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("D8").Value = "Jamarii"
$ws.Range("G8").Value = "Black or African American"

$ws.Range("B9").Value = 30
$ws.Range("C9").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("D9").Value = "Matthew"
$ws.Range("G9").Value = "White"

$ws.Range("B10").Value = 33
$ws.Range("C10").Value = "60b322994d0b901954690036"
$ws.Range("D10").Value = "Brennan"
$ws.Range("G10").Value = "White"
